$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change
$ws.Range("B3").Value = 17343333720216.06
$ws.Range("C3").Value = 20837433634322.24
$ws.Range("D3").Value = 3749834179261.602

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 15627226451574.46
$ws.Range("C4").Value = 4190871150888.41
$ws.Range("D4").Value = 2584225138292.988

# Row 5: AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 371090696275807.1
$ws.Range("C5").Value = 95342836577320.12
$ws.Range("D5").Value = 115422330140575.9
